$d = $word.ActiveDocument

$pairs = @(
    @("2025-05-09 Friday", "2025-05-10 Saturday"),
    @("51-23=", "14+21="),
    @("12+66=", "53+13="),
    @("21+56=", "21+13="),
    @("42-16=", "83-34="),
    @("85-23=", "53+35="),
    @("15+71=", "28-19="),
    @("53+0=", "11+45="),
    @("20+2=", "53+1="),
    @("96-1=", "53+19="),
    @("26+16=", "29-12="),
    @("28+16=", "48-12="),
    @("92-92=", "72+3="),
    @("43-39=", "37+48="),
    @("77+15=", "48-6="),
    @("39-38=", "3+78="),
    @("10+41=", "41-12="),
    @("10+52=", "69-32="),
    @("91-13=", "57+7="),
    @("64+17=", "79-74="),
    @("59+21=", "10+72="),
    @("31+25=", "65-17="),
    @("48-15=", "74+16="),
    @("8+17=", "90+8="),
    @("56-11=", "59+11="),
    @("41+7=", "66-28="),
    @("68+20=", "43+1="),
    @("37+61=", "7+4="),
    @("77+14=", "43+46="),
    @("83-51=", "33-17="),
    @("72-12=", "25+21="),
    @("2+63=", "28+2="),
    @("99-0=", "1+47="),
    @("7+62=", "35-10="),
    @("83+6=", "65-63="),
    @("45+46=", "57+20="),
    @("0+97=", "7+31="),
    @("21-2=", "24-21="),
    @("68-20=", "58+24="),
    @("52-50=", "3+48="),
    @("65-23=", "49-19="),
    @("82-19=", "28+63="),
    @("16+42=", "0+50="),
    @("61-21=", "62-57="),
    @("18+58=", "92-69="),
    @("3+16=", "22-3="),
    @("38+29=", "52-25="),
    @("28-13=", "55-17="),
    @("68-17=", "81-50="),
    @("77-22=", "21+73="),
    @("24+30=", "21+67="),
    @("67-1=", "98-5="),
    @("25+57=", "68-11="),
    @("69+2=", "69-9="),
    @("49-24=", "18+25="),
    @("59+4=", "93-92="),
    @("76+0=", "65+7="),
    @("31+58=", "41+47="),
    @("55-13=", "88-62="),
    @("49+32=", "34+57="),
    @("73-55=", "16+49="),
    @("13+39=", "30+57="),
    @("96-51=", "4+23="),
    @("39-19=", "53+34="),
    @("70-30=", "70-59="),
    @("81-61=", "55+41="),
    @("0+60=", "77-37="),
    @("60-25=", "38+49="),
    @("58-0=", "89-69="),
    @("87-85=", "29+26="),
    @("77-9=", "83+2="),
    @("53-23=", "48+16="),
    @("97-15=", "29-13="),
    @("43+33=", "11-3="),
    @("84+8=", "81-69="),
    @("30+42=", "99-9="),
    @("26+31=", "83-14="),
    @("84-45=", "65+9="),
    @("42-38=", "60-29="),
    @("92-75=", "61+35="),
    @("54-20=", "90-63="),
    @("57-6=", "39+24="),
    @("25+11=", "26+43="),
    @("98-27=", "31+49="),
    @("70-15=", "11+75="),
    @("15-6=", "40+20="),
    @("46-45=", "94-93="),
    @("30+20=", "63+10="),
    @("62-6=", "95-87="),
    @("48-13=", "20+67="),
    @("26+70=", "10+40="),
    @("56-14=", "95-87="),
    @("43-11=", "68-48="),
    @("16+41=", "74-67="),
    @("72-51=", "96-73="),
    @("65-3=", "6+29="),
    @("3+53=", "78-39="),
    @("50-37=", "65+12="),
    @("15+62=", "66+10="),
    @("45-41=", "76+23="),
    @("2+68=", "11+33="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"